$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.020.92'
$ws.Range('E2').Value = '  -0.51%  '

$ws.Range('D3').Value = '1.640.31'
$ws.Range('E3').Value = '  -0.90%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.18%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5052'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.28%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.009'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.20%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2576'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.36%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06447'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.22%  '

$ws.Range('E10').Value = '  -1.95%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07713'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.92%  '

$ws.Range('D12').Value = '1.648.83'
$ws.Range('E12').Value = '  -0.53%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.246'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.00%  '

$ws.Range('D14').Value = '1.865.37'
$ws.Range('E14').Value = '  -0.98%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5446'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.46%  '

$ws.Range('D16').Value = '0.0₅7922'
$ws.Range('E16').Value = '  -1.58%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.42'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.39%  '

$ws.Range('D18').Value = '26.009.59'
$ws.Range('E18').Value = '  -0.70%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.008'
$ws.Range('D19').ClearFormats()

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '202.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.89%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.283'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.19%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.991'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.70%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.967'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.15%  '

$ws.Range('E24').Value = '  +0.25%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.960'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +11.29%  '

$ws.Range('E26').Value = '  -1.54%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1150'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.99%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.75'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.09%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.718'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.56%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05051'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.18%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.240'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.00%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.251'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.20%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.195'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.89%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.539'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.14%  '

$ws.Range('E35').Value = '  -0.82%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.637'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.63%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.8899'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.59%  '

$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5619'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.63%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.148.31'
$ws.Range('E39').Value = '  -1.51%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01573'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.05%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.565'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.23%  '

$ws.Range('E42').Value = '  +0.25%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.669'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.19%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8077'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.46%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '99.67'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.27%  '

$ws.Range('D46').Value = '1.776.93'
$ws.Range('E46').Value = '  -0.97%  '

$ws.Range('E47').Value = '  +2.01%  '

$ws.Range('E48').Value = '  +0.46%  '

$ws.Range('E49').Value = '  +0.06%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '54.98'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.65%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05037'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.82%  '
